$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that used to sit right
#    after the title (Heading1) paragraph at the top of the document.
#    Range.Delete() removes the paragraph mark too, so the following
#    "Visuals and Sound Effects" heading simply slides up.
# ------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# ------------------------------------------------------------------
# 2. Add a new, bold "Play Book of Santa Slot Game for Free - Review"
#    paragraph right before the final (italic) paragraph, i.e. right
#    after the "No progressive jackpot feature" bullet.
#
#    To keep the new paragraph free of any inherited list/heading
#    style or italic formatting, it is first typed in a "clean" spot
#    (after a plain body paragraph with no direct character
#    formatting) and then moved into place with Cut/Paste - exactly
#    like a user would build it elsewhere and relocate it.
# ------------------------------------------------------------------
$stagingAnchor = $d.Paragraphs.Item(3)
$stagingAnchor.Range.InsertParagraphAfter()

$staging = $d.Paragraphs.Item($stagingAnchor.Index + 1)
$staging.Range.Text = "Play Book of Santa Slot Game for Free - Review"

$staging2 = $d.Paragraphs.Item($stagingAnchor.Index + 1)
$boldRange = $d.Range($staging2.Range.Start, $staging2.Range.End - 1)
$boldRange.Font.Bold = $true

$stagingFull = $d.Paragraphs.Item($stagingAnchor.Index + 1)
$stagingFull.Range.Cut()

$last = $d.Paragraphs.Last
$insertPoint = $d.Range($last.Range.Start, $last.Range.Start)
$insertPoint.Paste()

# ------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph - the old
#    image-generation prompt becomes the real meta-description text,
#    while the italic run formatting is left untouched.
# ------------------------------------------------------------------
$last2 = $d.Paragraphs.Last
$last2.Range.Find.Execute(
    "Create a feature image for Book of Santa that captures the essence of the game's theme and unique character. The image should be in a cartoon style and feature a happy-looking Maya warrior wearing glasses, as the protagonist of the game. The Maya warrior should be holding a large book in his hands, which should resemble Santa's book, with a few gifts spilling out of it. The background of the image should be a cozy fireplace scene, with the grid of the slot game superimposed on it. The image should use bright and cheerful colors and convey the festive mood of the holiday season.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Book of Santa slot game, play it for free and enter the Christmas atmosphere with expandable symbols and free spin feature.",
    2)
